# Generate Report for Handback
#
# This script applies the "handback" report-generation update to the
# localization-status workbook:
#   1. The global "Ready for handoff" status label is replaced everywhere
#      with "Handed back: in sync with en-US" (Overview sheet's status
#      columns plus each language sheet's Status column).
#   2. Each language sheet (zh-cn, de-de) gains two new populated columns:
#        F = "Latest Target File"   (same file reference as column A)
#        G = "Latest Handback File" (same file reference as column D)
#      formatted/linked the same way as the existing hyperlink columns.
#   3. The "Latest Handback DateTime" column (H) on each language sheet is
#      updated from the "never handed back" placeholder to the actual
#      handback timestamp for that language.

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# --- 1. Update the status label everywhere it is used -----------------
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- helper values reused for the new F/G columns ----------------------
$hdrFirstMd   = "5b561573-67ca-4fef-be14-49acd8225084.md"
$hdrSecondMd  = "ef958983-997d-40ee-8cde-e7336ffd1a92.md"

$zhCnFirstXlf  = "5b561573-67ca-4fef-be14-49acd8225084.44e30c1cc3d27bb12124d6f9e60edc59c9c6f778.zh-cn.xlf"
$zhCnSecondXlf = "ef958983-997d-40ee-8cde-e7336ffd1a92.d488aa41eb9bd1dbf8e827f3c2ebbdebc604a0b2.zh-cn.xlf"
$deDeFirstXlf  = "5b561573-67ca-4fef-be14-49acd8225084.44e30c1cc3d27bb12124d6f9e60edc59c9c6f778.de-de.xlf"
$deDeSecondXlf = "ef958983-997d-40ee-8cde-e7336ffd1a92.d488aa41eb9bd1dbf8e827f3c2ebbdebc604a0b2.de-de.xlf"

$srcRepoFirstUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/a3e23e85da0d9bdeae55b9cf20c2a0c1a40e0e6a/e2e/5b561573-67ca-4fef-be14-49acd8225084.md"
$srcRepoSecondUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a3e23e85da0d9bdeae55b9cf20c2a0c1a40e0e6a/e2e/ef958983-997d-40ee-8cde-e7336ffd1a92.md"
$zhCnFirstUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ce91609be29abf9e5404fb69c62677b2940cd90/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/5b561573-67ca-4fef-be14-49acd8225084.44e30c1cc3d27bb12124d6f9e60edc59c9c6f778.zh-cn.xlf"
$zhCnSecondUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ce91609be29abf9e5404fb69c62677b2940cd90/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/ef958983-997d-40ee-8cde-e7336ffd1a92.d488aa41eb9bd1dbf8e827f3c2ebbdebc604a0b2.zh-cn.xlf"
$deDeFirstUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c4475c5d267399990d7fc912b9d7bb48426797e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/5b561573-67ca-4fef-be14-49acd8225084.44e30c1cc3d27bb12124d6f9e60edc59c9c6f778.de-de.xlf"
$deDeSecondUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c4475c5d267399990d7fc912b9d7bb48426797e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/ef958983-997d-40ee-8cde-e7336ffd1a92.d488aa41eb9bd1dbf8e827f3c2ebbdebc604a0b2.de-de.xlf"

$hyperlinkUnderline = 2
$hyperlinkColor = 15570276

function Set-HandbackLink($ws, $cellRef, $address, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText)
    $ws.Range($cellRef).Font.Underline = $hyperlinkUnderline
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# --- 2. Populate the new "Latest Target File" / "Latest Handback File" --
#        columns (F/G) on each language sheet ---------------------------
Set-HandbackLink $wsZhCn "F2" $srcRepoFirstUrl  $hdrFirstMd
Set-HandbackLink $wsZhCn "G2" $zhCnFirstUrl     $zhCnFirstXlf
Set-HandbackLink $wsZhCn "F3" $srcRepoSecondUrl $hdrSecondMd
Set-HandbackLink $wsZhCn "G3" $zhCnSecondUrl    $zhCnSecondXlf

Set-HandbackLink $wsDeDe "F2" $srcRepoFirstUrl  $hdrFirstMd
Set-HandbackLink $wsDeDe "G2" $deDeFirstUrl     $deDeFirstXlf
Set-HandbackLink $wsDeDe "F3" $srcRepoSecondUrl $hdrSecondMd
Set-HandbackLink $wsDeDe "G3" $deDeSecondUrl    $deDeSecondXlf

# --- 3. Record the actual handback timestamps on column H --------------
$wsZhCn.Range("H2").Value = "2016-03-11 08:59:40"
$wsZhCn.Range("H3").Value = "2016-03-11 08:59:40"

$wsDeDe.Range("H2").Value = "2016-03-11 08:59:46"
$wsDeDe.Range("H3").Value = "2016-03-11 08:59:46"

Write-Host "Handback report generated."
